# Generate Report for Handoff
#
# Source file "81a5b25e-1357-4f92-ab21-31b8813c4143.md" moved from
# "In Translation" to "Ready for handoff" for both locales, and the
# localization-status report timestamps (and the report columns that
# display them) were refreshed to reflect the new handoff generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Status columns for zh-cn (E) and de-de (F) move to "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"

# Latest HO Xliff Generate Date (G) advances to the new handoff run time
$wsOverview.Range("G2").Value = "2016-08-30 18:49:02"

# Widen the two status columns to fit the new, longer status text
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-30 18:48:56"
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-30 18:49:02"
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
